# Ion trace tables (tabs 1 & 2: ions_rain_LD / ions_rain_LL) dropped the
# leading "rowname" index column and the "Time" reference row, shifting the
# remaining Ion/pVal/phase/peak.shape/period data up-and-left by one.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Delete()      # remove the "Time" data row
$ws1.Columns.Item(1).Delete()   # remove the "rowname" index column

$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Delete()      # remove the "Time" data row
$ws2.Columns.Item(1).Delete()   # remove the "rowname" index column
